# summary_of_experiments.xlsx - "Add files via upload" edit
#
# The uploaded revision rewrote the description of experiment #4's model
# (cell D5) from a single "A Neural Network" to an "Ensemble of 5 Neural
# Networks" with the same architecture bullet list, which nudged rows 4/5
# to a custom height and moved the active selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newModelText = @"
Ensemble of 5 Neural Networks with architecture:
1. Embedding Layer
2. BiLSTM Layer
3. GlobalMaxPooling1D 
4. BatchNormalization
5. Drop out of 0.5
6. Dense Layer with relu activation
7. Drop out of 0.5
8. Dense with with relu activation
9. Dropout of 0.5
10. Dense with sigmoid activation
"@

$ws.Range("D5").Value = $newModelText

# Match the left/top-wrapped alignment already used by the other long
# text cells in this table (e.g. C5/B6/C6/E6).
$ws.Range("D5").VerticalAlignment = -4160   # xlTop
$ws.Range("D5").HorizontalAlignment = -4131 # xlLeft
$ws.Range("D5").WrapText = $true

# The longer row 4 text now fits in a shorter custom row height, and row 5
# grew slightly to fit the new wording.
$ws.Rows.Item(4).RowHeight = 144.75
$ws.Rows.Item(5).RowHeight = 142.5

# Selection ended up on D6 (and scrolled back so row 1 is visible again).
$ws.Range("D6").Select()
